# Yellow highlighted validated cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update / recompute the validated values in column C ---
$ws.Range("C2").Formula = "=(44.9 + 40.3) / 2"
$ws.Range("C5").Value = 46
$ws.Range("C10").Value = 4.76
$ws.Range("C11").Value = 4.28
$ws.Range("C14").Value = 1.8
$ws.Range("C19").Value = 5
$ws.Range("C23").Value = 7.5
$ws.Range("C24").Formula = "=46*0.48"
$ws.Range("C27").Value = 78400
$ws.Range("C28").Value = 22.21

# --- Highlight all validated cells in column C with a solid yellow fill ---
$validatedCells = @("C2","C5","C10","C11","C13","C14","C15","C16","C18","C19","C22","C23","C24","C27","C28")
foreach ($addr in $validatedCells) {
    $ws.Range($addr).Interior.Color = 65535
}

# --- Restore the sheet view: no frozen scroll offset, selection on C10 ---
$null = $ws.Range("C10").Select()
